# "finalização plano iteração elaboração"
#
# The Elaboração iteration's two remaining "E2" work items — analysis/design
# of the use case, and the resulting use-case design itself — are wrapped
# up: their status moves from "Iniciado" to "Finalizado" and the hours
# actually worked are updated to reflect the extra effort spent finishing
# them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 - "Analise e projeto do caso de uso": mark finished, log worked hours.
$ws.Range("D20").Value = "Finalizado"
$ws.Range("H20").Value = 10

# Row 21 - "Caso de uso projetado": mark finished, log worked hours.
$ws.Range("D21").Value = "Finalizado"
$ws.Range("H21").Value = 32

# The filtered range is re-confirmed as part of closing out the iteration,
# which stamps a fresh (local) _FilterDatabase name alongside the existing
# one.
$ws.Names.Add("_xlnm._FilterDatabase_0", "='Lista de Itens de Trabalho'!`$A`$1:`$I`$25")

# Move the view/selection down onto the row below the now-finished items.
$ws.Range("H22").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
